$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 598 held a stray duplicate of "KSB" (it already exists correctly at
# row 390); the author overwrote it with the "Purvan" ticker that used to
# be the last row of the list.
$ws.Range("A598").Value = "Purvan"

# Six brand-new F&O tickers appended right after it.
$ws.Range("A599").Value = "PGEL"
$ws.Range("A600").Value = "NIFTY"
$ws.Range("A601").Value = "BANKNIFTY"
$ws.Range("A602").Value = "IREDA"
$ws.Range("A603").Value = "GMRAIRPORT"
$ws.Range("A604").Value = "ETERNAL"

# Those six cells were blank placeholders styled with s="1"; typing fresh
# values into them in Excel drops that inherited formatting.
$ws.Range("A599:A604").Style = "Normal"

# The rest of the old styled-but-empty placeholder cells (605-625) are
# cleared out completely, becoming plain blank rows.
$ws.Range("A605:A625").Clear()

# The sheet's last row (1000) is removed, shrinking the used range to A1:A999.
$ws.Rows.Item(1000).Delete()

# Restore the on-screen selection/scroll position to match the saved view.
$ws.Range("C607").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 576
$win.ScrollColumn = 1
